$wb = $excel.ActiveWorkbook

# Remove the "2020-10-07" service entry (id 3) from sheet "test 7418"
$ws7418 = $wb.Worksheets.Item("test 7418")
$ws7418.Rows(2).Delete()
$ws7418.Columns(3).AutoFit()

# "test 7418" now appears before "test 640" in the tab order
$ws7418.Move($wb.Worksheets.Item(1))
